$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '27.999.89'
$ws.Cells.Item(2, 5).Value = '  -3.76%  '

$ws.Cells.Item(3, 4).Value = '1.868.21'
$ws.Cells.Item(3, 5).Value = '  -2.90%  '

$ws.Cells.Item(4, 4).Value = '1.004'
$ws.Cells.Item(4, 5).Value = '  -0.32%  '

$ws.Cells.Item(5, 4).Value = '318.77'
$ws.Cells.Item(5, 5).Value = '  -2.28%  '

$ws.Cells.Item(6, 4).Value = '1.004'
$ws.Cells.Item(6, 5).Value = '  -0.30%  '

$ws.Cells.Item(7, 5).Value = '  -5.88%  '

$ws.Cells.Item(8, 4).Value = '0.3736'
$ws.Cells.Item(8, 5).Value = '  -2.35%  '

$ws.Cells.Item(9, 4).Value = '0.07446'
$ws.Cells.Item(9, 5).Value = '  -4.10%  '

$ws.Cells.Item(10, 4).Value = '0.9320'
$ws.Cells.Item(10, 5).Value = '  -4.72%  '

$ws.Cells.Item(11, 4).Value = '21.22'
$ws.Cells.Item(11, 5).Value = '  -6.64%  '

$ws.Cells.Item(12, 4).Value = '1.964.98'
$ws.Cells.Item(12, 5).Value = '  +3.12%  '

$ws.Cells.Item(13, 4).Value = '6.730'
$ws.Cells.Item(13, 5).Value = '  -3.46%  '

$ws.Cells.Item(14, 4).Value = '5.433'
$ws.Cells.Item(14, 5).Value = '  -4.60%  '

$ws.Cells.Item(15, 4).Value = '0.06867'
$ws.Cells.Item(15, 5).Value = '  -2.84%  '

$ws.Cells.Item(16, 4).Value = '1.004'
$ws.Cells.Item(16, 5).Value = '  -0.45%  '

$ws.Cells.Item(17, 4).Value = '80.89'
$ws.Cells.Item(17, 5).Value = '  -4.10%  '

$ws.Cells.Item(18, 4).Value = '0.000009027'
$ws.Cells.Item(18, 5).Value = '  -5.13%  '

$ws.Cells.Item(19, 4).Value = '1.004'
$ws.Cells.Item(19, 5).Value = '  -0.26%  '

$ws.Cells.Item(20, 5).Value = '  -6.06%  '

$ws.Cells.Item(21, 4).Value = '27.987.67'
$ws.Cells.Item(21, 5).Value = '  -3.86%  '

$ws.Cells.Item(22, 4).Value = '5.123'

$ws.Cells.Item(23, 4).Value = '10.99'
$ws.Cells.Item(23, 5).Value = '  +0.13%  '

$ws.Cells.Item(24, 4).Value = '2.134.57'
$ws.Cells.Item(24, 5).Value = '  +0.16%  '

$ws.Cells.Item(25, 4).Value = '2.039'
$ws.Cells.Item(25, 5).Value = '  -1.56%  '

$ws.Cells.Item(26, 4).Value = '153.62'
$ws.Cells.Item(26, 5).Value = '  -2.83%  '

$ws.Cells.Item(27, 4).Value = '18.52'
$ws.Cells.Item(27, 5).Value = '  -3.14%  '

$ws.Cells.Item(28, 4).Value = '5.509'
$ws.Cells.Item(28, 5).Value = '  -2.66%  '

$ws.Cells.Item(29, 4).Value = '113.28'
$ws.Cells.Item(29, 5).Value = '  -4.12%  '

$ws.Cells.Item(30, 4).Value = '1.697'
$ws.Cells.Item(30, 5).Value = '  -7.53%  '

$ws.Cells.Item(31, 4).Value = '0.08989'
$ws.Cells.Item(31, 5).Value = '  -3.78%  '

$ws.Cells.Item(32, 4).Value = '0.8077'
$ws.Cells.Item(32, 5).Value = '  -5.47%  '

$ws.Cells.Item(33, 4).Value = '4.788'
$ws.Cells.Item(33, 5).Value = '  -6.44%  '

$ws.Cells.Item(34, 4).Value = '1.177'
$ws.Cells.Item(34, 5).Value = '  -5.36%  '

$ws.Cells.Item(35, 4).Value = '2.953'
$ws.Cells.Item(35, 5).Value = '  -2.35%  '

$ws.Cells.Item(36, 5).Value = '  -0.22%  '

$ws.Cells.Item(37, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(37, 4).Value = '1.120'
$ws.Cells.Item(37, 5).Value = '  -3.61%  '

$ws.Cells.Item(38, 2).Value = 'Hedera'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(38, 4).Value = '0.05492'
$ws.Cells.Item(38, 5).Value = '  -3.38%  '

$ws.Cells.Item(39, 2).Value = 'VeChain'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(39, 4).Value = '0.01972'
$ws.Cells.Item(39, 5).Value = '  -3.62%  '

$ws.Cells.Item(40, 2).Value = 'MXToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(40, 4).Value = '2.997'
$ws.Cells.Item(40, 5).Value = '  -5.34%  '

$ws.Cells.Item(41, 4).Value = '0.5244'
$ws.Cells.Item(41, 5).Value = '  -4.96%  '

$ws.Cells.Item(42, 4).Value = '6.997'
$ws.Cells.Item(42, 5).Value = '  -6.74%  '

$ws.Cells.Item(43, 4).Value = '0.1686'
$ws.Cells.Item(43, 5).Value = '  -4.08%  '

$ws.Cells.Item(44, 4).Value = '8.771'
$ws.Cells.Item(44, 5).Value = '  -5.97%  '

$ws.Cells.Item(45, 4).Value = '0.06722'
$ws.Cells.Item(45, 5).Value = '  -2.98%  '

$ws.Cells.Item(46, 4).Value = '0.4877'
$ws.Cells.Item(46, 5).Value = '  -6.25%  '

$ws.Cells.Item(47, 4).Value = '10.52'
$ws.Cells.Item(47, 5).Value = '  -6.18%  '

$ws.Cells.Item(48, 4).Value = '106.81'
$ws.Cells.Item(48, 5).Value = '  -3.26%  '

$ws.Cells.Item(49, 5).Value = '  -0.41%  '

$ws.Cells.Item(50, 4).Value = '1.673'
$ws.Cells.Item(50, 5).Value = '  -5.43%  '

$ws.Cells.Item(51, 4).Value = '1.876'
$ws.Cells.Item(51, 5).Value = '  -14.28%  '
